$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.046.35'
$ws.Range("E2").Value = '  -0.89%  '
$ws.Range("D3").Value = '1.824.80'
$ws.Range("E3").Value = '  -0.01%  '
$ws.Range("E4").Value = '  -0.47%  '
$ws.Range("D5").Value = '''311.90'
$ws.Range("E5").Value = '  -1.01%  '
$ws.Range("E6").Value = '  -0.36%  '
$ws.Range("D7").Value = '''0.4403'
$ws.Range("E7").Value = '  +2.64%  '
$ws.Range("D8").Value = '''0.3682'
$ws.Range("E8").Value = '  -0.23%  '
$ws.Range("D9").Value = '''0.07278'
$ws.Range("E9").Value = '  +0.55%  '
$ws.Range("D10").Value = '''0.8470'
$ws.Range("E10").Value = '  -1.91%  '
$ws.Range("D11").Value = '''20.70'
$ws.Range("E11").Value = '  -1.89%  '
$ws.Range("D12").Value = '1.821.35'
$ws.Range("E12").Value = '  -0.17%  '
$ws.Range("D13").Value = '''6.666'
$ws.Range("E13").Value = '  -0.01%  '
$ws.Range("D14").Value = '''0.07074'
$ws.Range("E14").Value = '  +0.28%  '
$ws.Range("D15").Value = '''5.305'
$ws.Range("E15").Value = '  -0.72%  '
$ws.Range("D16").Value = '''90.21'
$ws.Range("E16").Value = '  +3.09%  '
$ws.Range("D17").Value = '''1.001'
$ws.Range("E17").Value = '  -0.50%  '
$ws.Range("D18").Value = '''0.000008801'
$ws.Range("E18").Value = '  -0.74%  '
$ws.Range("E19").Value = '  -0.33%  '
$ws.Range("E20").Value = '  -1.64%  '
$ws.Range("D21").Value = '27.107.11'
$ws.Range("E21").Value = '  -0.78%  '
$ws.Range("D22").Value = '''5.153'
$ws.Range("E22").Value = '  +0.00%  '
$ws.Range("E23").Value = '  +0.33%  '
$ws.Range("D24").Value = '2.047.33'
$ws.Range("E24").Value = '  -0.25%  '
$ws.Range("D25").Value = '''1.998'
$ws.Range("E25").Value = '  -0.58%  '
$ws.Range("D26").Value = '''151.65'
$ws.Range("E26").Value = '  -0.79%  '
$ws.Range("E27").Value = '  +3.02%  '
$ws.Range("D28").Value = '''18.28'
$ws.Range("E28").Value = '  -0.56%  '
$ws.Range("D29").Value = '''5.235'
$ws.Range("E29").Value = '  -0.53%  '
$ws.Range("D30").Value = '''117.05'
$ws.Range("E30").Value = '  +0.30%  '
$ws.Range("E31").Value = '  -0.59%  '
$ws.Range("E32").Value = '  -1.49%  '
$ws.Range("D33").Value = '''0.7416'
$ws.Range("E33").Value = '  -2.90%  '
$ws.Range("D34").Value = '''4.428'
$ws.Range("E34").Value = '  -1.23%  '
$ws.Range("E35").Value = '  +1.13%  '
$ws.Range("E36").Value = '  -0.44%  '
$ws.Range("E37").Value = '  -2.07%  '
$ws.Range("D38").Value = '''0.01950'
$ws.Range("E38").Value = '  -0.26%  '
$ws.Range("D39").Value = '''0.05239'
$ws.Range("E39").Value = '  -0.36%  '
$ws.Range("D40").Value = '''7.272'
$ws.Range("E40").Value = '  +2.13%  '
$ws.Range("D41").Value = '''2.869'
$ws.Range("E41").Value = '  -0.20%  '
$ws.Range("D42").Value = '''0.5173'
$ws.Range("E42").Value = '  +2.36%  '
$ws.Range("E43").Value = '  +1.22%  '
$ws.Range("D44").Value = '''8.542'
$ws.Range("E44").Value = '  -1.17%  '
$ws.Range("D45").Value = '''10.61'
$ws.Range("E45").Value = '  +0.94%  '
$ws.Range("D46").Value = '''0.4828'
$ws.Range("E46").Value = '  +2.48%  '
$ws.Range("D47").Value = '''106.08'
$ws.Range("E47").Value = '  +0.18%  '
$ws.Range("D48").Value = '''1.934'
$ws.Range("E48").Value = '  +6.64%  '
$ws.Range("D49").Value = '''0.9998'
$ws.Range("E49").Value = '  -0.38%  '
$ws.Range("D50").Value = '''0.06335'
$ws.Range("E50").Value = '  -1.28%  '
$ws.Range("D51").Value = '''1.660'
$ws.Range("E51").Value = '  -0.29%  '
